$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append, continuing the daily series after row 244 (2021-05-02).
$newRows = @(
    @{ Row = 245; A = 44319; B = 0; C = 5; D = 127.1617497456765 },
    @{ Row = 246; A = 44320; B = 0; C = 4; D = 101.7293997965412 },
    @{ Row = 247; A = 44321; B = 0; C = 4; D = 101.7293997965412 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Column A keeps the same style/number-format as the rows above it
    # (copy the formatting down before writing the value).
    $ws.Range("A244").Copy() | Out-Null
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
}

$excel.CutCopyMode = 0
